$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add a new row of data (row 11) continuing the time-tracking log.
# Copy the date formatting used by the cell above so we reuse the same style.
$ws.Range("A10").Copy($ws.Range("A11"))
$ws.Range("A11").Value = 42732

$ws.Range("B11").Value = "8h"
$ws.Range("C11").Value = "Kleinere Anpassungen an der Oberfläche, Anzeigen der Kurse hinzugefügt, refactoring und Funktionalitäten im ViewModelbereich"

# Update the active selection to reflect the next empty row, like Excel would.
$ws.Range("C12").Select()
